$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 33.727272
$ws.Range("I11").Value = 33.727272
$ws.Range("K11").Value = 33.727272
$ws.Range("M11").Value = 106.272728

$ws.Range("H38").Value = 15175.5
$ws.Range("I38").Value = 19999.334
$ws.Range("J38").Value = 704
$ws.Range("K38").Value = 59998.00199999999
$ws.Range("L38").Value = 2112
$ws.Range("M38").Value = -59626.00199999999
$ws.Range("N38").Value = -2856

$ws.Range("H41").Value = 681.2
$ws.Range("I41").Value = 316.2857
$ws.Range("J41").Value = 1000.5
$ws.Range("K41").Value = 316.2857
$ws.Range("L41").Value = 1000.5
$ws.Range("M41").Value = 123.7143
$ws.Range("N41").Value = -1880.5

$ws.Range("H116").Value = 11637.477
$ws.Range("I116").Value = 10213.286
$ws.Range("J116").Value = 14485.857
$ws.Range("K116").Value = 10213.286
$ws.Range("L116").Value = 14485.857
$ws.Range("M116").Value = -6771.286
$ws.Range("N116").Value = -21369.857

$ws.Range("H132").Value = 1394.0364
$ws.Range("I132").Value = 1099.3529
$ws.Range("K132").Value = 3298.0587
$ws.Range("M132").Value = -768.0587000000005

$ws.Range("H137").Value = 2547.7932
$ws.Range("I137").Value = 2176.3333
$ws.Range("K137").Value = 6528.999899999999
$ws.Range("M137").Value = -3978.999899999999

$ws.Range("H138").Value = 1282405.8
$ws.Range("J138").Value = 1485799.8
$ws.Range("L138").Value = 4457399.4
$ws.Range("N138").Value = -4467679.4

$ws.Range("H141").Value = 7372.1875
$ws.Range("I141").Value = 4164.9165
$ws.Range("K141").Value = 12494.7495
$ws.Range("M141").Value = -7314.749500000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 71502.984
$ws.Range("I32").Value = 72886.63
$ws.Range("K32").Value = 72886.63
$ws.Range("M32").Value = -72599.63

$ws.Range("H45").Value = 10628.385
$ws.Range("I45").Value = 12116.9
$ws.Range("J45").Value = 5666.6665
$ws.Range("K45").Value = 12116.9
$ws.Range("L45").Value = 5666.6665
$ws.Range("M45").Value = -11739.9
$ws.Range("N45").Value = -6420.6665

$ws.Range("H52").Value = 85000
$ws.Range("J52").Value = 85000
$ws.Range("L52").Value = 85000
$ws.Range("N52").Value = -85636

$ws.Range("H132").Value = 10671.846
$ws.Range("J132").Value = 5038
$ws.Range("L132").Value = 15114
$ws.Range("N132").Value = -20174

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7503.8184
$ws.Range("I86").Value = 11071.5
$ws.Range("J86").Value = 3222.6
$ws.Range("K86").Value = 11071.5
$ws.Range("L86").Value = 3222.6
$ws.Range("M86").Value = -9948.5
$ws.Range("N86").Value = -5468.6

$ws.Range("H89").Value = 7503.8184
$ws.Range("I89").Value = 11071.5
$ws.Range("J89").Value = 3222.6
$ws.Range("K89").Value = 55357.5
$ws.Range("L89").Value = 16113
$ws.Range("M89").Value = -49741.5
$ws.Range("N89").Value = -27345

$ws.Range("H94").Value = 2128.4285
$ws.Range("I94").Value = 2128.4285
$ws.Range("K94").Value = 2128.4285
$ws.Range("M94").Value = -1677.4285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4246.643
$ws.Range("I31").Value = 2540.5
$ws.Range("J31").Value = 7019.125
$ws.Range("K31").Value = 2540.5
$ws.Range("L31").Value = 7019.125
$ws.Range("M31").Value = -2245.5
$ws.Range("N31").Value = -7609.125

$ws.Range("H34").Value = 4246.643
$ws.Range("I34").Value = 2540.5
$ws.Range("J34").Value = 7019.125
$ws.Range("K34").Value = 2540.5
$ws.Range("L34").Value = 7019.125
$ws.Range("M34").Value = -2338.5
$ws.Range("N34").Value = -7423.125

$ws.Range("H58").Value = 2296.8386
$ws.Range("I58").Value = 2363.9333
$ws.Range("K58").Value = 2363.9333
$ws.Range("M58").Value = -2160.9333

$ws.Range("H62").Value = 15000
$ws.Range("I62").Value = 15000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 15000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -14376
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 15000
$ws.Range("I65").Value = 15000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 75000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -71880
$ws.Range("N65").Value = -71880

$ws.Range("H136").Value = 2296.8386
$ws.Range("I136").Value = 2363.9333
$ws.Range("K136").Value = 7091.7999
$ws.Range("M136").Value = -4541.7999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 350.5
$ws.Range("J25").Value = 501
$ws.Range("L25").Value = 1503
$ws.Range("N25").Value = -1841

$ws.Range("H30").Value = 350.5
$ws.Range("J30").Value = 501
$ws.Range("L30").Value = 1503
$ws.Range("N30").Value = -1707

$ws.Range("H131").Value = 5871.48
$ws.Range("I131").Value = 1138.5714
$ws.Range("J131").Value = 7712.0557
$ws.Range("K131").Value = 3415.7142
$ws.Range("L131").Value = 23136.1671
$ws.Range("M131").Value = 1624.2858
$ws.Range("N131").Value = -33216.1671

$ws.Range("H140").Value = 1749.4642
$ws.Range("I140").Value = 1699.5
$ws.Range("K140").Value = 5098.5
$ws.Range("M140").Value = 81.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 500027500
$ws.Range("J34").Value = 500027500
$ws.Range("L34").Value = 500027500
$ws.Range("N34").Value = -500028036

$ws.Range("H76").Value = 500027500
$ws.Range("J76").Value = 500027500
$ws.Range("L76").Value = 500027500
$ws.Range("N76").Value = -500028130

$ws.Range("H79").Value = 500027500
$ws.Range("J79").Value = 500027500
$ws.Range("L79").Value = 500027500
$ws.Range("N79").Value = -500029684

$ws.Range("H132").Value = 3216.7058
$ws.Range("I132").Value = 3329.5
$ws.Range("J132").Value = 3169.7083
$ws.Range("K132").Value = 9988.5
$ws.Range("L132").Value = 9509.124899999999
$ws.Range("M132").Value = -7458.5
$ws.Range("N132").Value = -14569.1249

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1075.5
$ws.Range("J16").Value = 498
$ws.Range("L16").Value = 498
$ws.Range("N16").Value = -838

$ws.Range("H22").Value = 3499
$ws.Range("I22").Value = 1248.5
$ws.Range("J22").Value = 5749.5
$ws.Range("K22").Value = 1248.5
$ws.Range("L22").Value = 5749.5
$ws.Range("M22").Value = -953.5
$ws.Range("N22").Value = -6339.5

$ws.Range("H27").Value = 3499
$ws.Range("I27").Value = 1248.5
$ws.Range("J27").Value = 5749.5
$ws.Range("K27").Value = 1248.5
$ws.Range("L27").Value = 5749.5
$ws.Range("M27").Value = -1141.5
$ws.Range("N27").Value = -5963.5

$ws.Range("H46").Value = 3297
$ws.Range("I46").Value = 2139.111
$ws.Range("K46").Value = 2139.111
$ws.Range("M46").Value = -1951.111

$ws.Range("H132").Value = 1217029.6
$ws.Range("I132").Value = 3036.8667
$ws.Range("K132").Value = 9110.6001
$ws.Range("M132").Value = -6580.6001

$ws.Range("H133").Value = 88579
$ws.Range("J133").Value = 88579
$ws.Range("L133").Value = 88579
$ws.Range("N133").Value = -93639

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12375
$ws.Range("I62").Value = 7000
$ws.Range("K62").Value = 7000
$ws.Range("M62").Value = -6376

$ws.Range("H65").Value = 12375
$ws.Range("I65").Value = 7000
$ws.Range("K65").Value = 35000
$ws.Range("M65").Value = -31880

$ws.Range("H126").Value = 2941.5
$ws.Range("I126").Value = 3029.8
$ws.Range("K126").Value = 9089.400000000001
$ws.Range("M126").Value = -6619.400000000001

$ws.Range("H132").Value = 6603.932
$ws.Range("I132").Value = 6343.7236
$ws.Range("K132").Value = 19031.1708
$ws.Range("M132").Value = -16501.1708

$ws.Range("H136").Value = 4528.6294
$ws.Range("I136").Value = 1638.7
$ws.Range("K136").Value = 4916.1
$ws.Range("M136").Value = -2366.1
